$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C4").Value = 8412
$ws.Range("C5:C6").Value = 8387
$ws.Range("C7:C13").Value = 8112
$ws.Range("C14:C31").Value = 7707
$ws.Range("C32:C63").Value = 7657
$ws.Range("C64:C252").Value = 7573
